$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.223.09'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.913.85'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8204'
$ws.Range('E5').Value = '  +3.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.04'
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3257'
$ws.Range('E8').Value = '  +3.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.92'
$ws.Range('E9').Value = '  +2.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07081'
$ws.Range('E10').Value = '  +2.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08095'
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7808'
$ws.Range('E12').Value = '  +4.68%  '
$ws.Range('D13').Value = '1.900.33'
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.311'
$ws.Range('E14').Value = '  +1.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.44'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').Value = '30.211.56'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.28'
$ws.Range('E17').Value = '  +1.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.945'
$ws.Range('E18').Value = '  +0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '247.30'
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007810'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').Value = '2.165.20'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.294'
$ws.Range('E24').Value = '  +5.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1686'
$ws.Range('E25').Value = '  +22.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.359'
$ws.Range('E26').Value = '  +0.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.61'
$ws.Range('E27').Value = '  -1.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.02'
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.121'
$ws.Range('E29').Value = '  +4.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.372'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.321'
$ws.Range('E32').Value = '  -0.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05773'
$ws.Range('E33').Value = '  +5.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.106'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.281'
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7437'
$ws.Range('E36').Value = '  +1.16%  '
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.713'
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01932'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4478'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.60'
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.978'
$ws.Range('E43').Value = '  -3.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8524'
$ws.Range('E44').Value = '  +1.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.915'
$ws.Range('E45').Value = '  +1.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9998'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.04'
$ws.Range('E47').Value = '  +2.54%  '
$ws.Range('D48').Value = '1.021.58'
$ws.Range('E48').Value = '  +4.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.611'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.857'
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.571'
$ws.Range('E51').Value = '  +5.00%  '
